$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: zero out B5:J5, keep K5 as 0, add L5 = 0
$ws.Range("B5:L5").Value = 0.0

# Row 6: update B6:K6 with new values, add L6 = 100100
$ws.Range("B6").Value = 15.0
$ws.Range("C6").Value = 10.0
$ws.Range("D6").Value = 1500.0
$ws.Range("E6").Value = 100000.0
$ws.Range("F6").Value = 10.0
$ws.Range("G6").Value = 10.0
$ws.Range("H6").Value = 10.0
$ws.Range("I6").Value = 80.0
$ws.Range("J6").Value = 0.0
$ws.Range("K6").Value = 0.0
$ws.Range("L6").Value = 100100.0
